# Auto update Excel log
# Appends newly logged sensor events to the PIR and Proximity sheets.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append bathroom "No Motion" readings, then a final
# "Motion Detected" reading, as rows 49-61 ---
$pir = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-02-06", "09:40:34", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:40:35", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:40:40", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:40:45", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:40:50", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:40:55", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:00", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:05", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:10", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:15", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:20", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:25", "09:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "09:41:28", "09:00", "Bathroom", "Motion Detected", "Active")
)

$startRow = 49
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $row = $pirRows[$i]

    # Column A holds a plain "YYYY-MM-DD" text value (not a real date), so
    # force text formatting to stop Excel from auto-converting it to a
    # date serial, then restore the Normal style so no stray per-cell
    # number format lingers on the written-out cell.
    $pir.Cells.Item($r, 1).NumberFormat = "@"
    $pir.Cells.Item($r, 1).Value = $row[0]
    $pir.Cells.Item($r, 1).Style = "Normal"

    $pir.Cells.Item($r, 2).Value = $row[1]
    $pir.Cells.Item($r, 3).Value = $row[2]
    $pir.Cells.Item($r, 4).Value = $row[3]
    $pir.Cells.Item($r, 5).Value = $row[4]
    $pir.Cells.Item($r, 6).Value = $row[5]
}

# --- Proximity sheet: append Bathroom Door ENTER event as row 5 ---
$prox = $wb.Worksheets.Item("Proximity")

$prox.Cells.Item(5, 1).NumberFormat = "@"
$prox.Cells.Item(5, 1).Value = "2026-02-06"
$prox.Cells.Item(5, 1).Style = "Normal"

$prox.Cells.Item(5, 2).Value = "09:41:34"
$prox.Cells.Item(5, 3).Value = "09:00"
$prox.Cells.Item(5, 4).Value = "Bathroom Door"
$prox.Cells.Item(5, 5).Value = "ENTER"
$prox.Cells.Item(5, 6).Value = "User ENTERED Bathroom"
